$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035981583171684
$ws.Cells.Item(2, 4).Value = 1.043259295599308
$ws.Cells.Item(2, 5).Value = 1.034964889398478
$ws.Cells.Item(2, 6).Value = 1.050552412143832
$ws.Cells.Item(2, 9).Value = 1.04062089495411
$ws.Cells.Item(2, 10).Value = 1.041092390456456
$ws.Cells.Item(2, 11).Value = 1.046033665589698
$ws.Cells.Item(2, 12).Value = 1.037762851298454
$ws.Cells.Item(2, 13).Value = 1.053306371980609
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036848356035098
$ws.Cells.Item(3, 4).Value = 1.043943365355207
$ws.Cells.Item(3, 5).Value = 1.035699518811473
$ws.Cells.Item(3, 6).Value = 1.051417456946366
$ws.Cells.Item(3, 9).Value = 1.04085231878188
$ws.Cells.Item(3, 10).Value = 1.041603354871949
$ws.Cells.Item(3, 11).Value = 1.046529066063775
$ws.Cells.Item(3, 12).Value = 1.038306976380272
$ws.Cells.Item(3, 13).Value = 1.053983750560194
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037409729516642
$ws.Cells.Item(4, 4).Value = 1.044386415440683
$ws.Cells.Item(4, 5).Value = 1.036175693641937
$ws.Cells.Item(4, 6).Value = 1.051978085168056
$ws.Cells.Item(4, 9).Value = 1.041000964028836
$ws.Cells.Item(4, 10).Value = 1.041933844084614
$ws.Cells.Item(4, 11).Value = 1.046849355131162
$ws.Cells.Item(4, 12).Value = 1.038659214109822
$ws.Cells.Item(4, 13).Value = 1.054422303318093
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037645852393097
$ws.Cells.Item(5, 4).Value = 1.0445727704747
$ws.Cells.Item(5, 5).Value = 1.036376072272099
$ws.Cells.Item(5, 6).Value = 1.052213983691402
$ws.Cells.Item(5, 9).Value = 1.04106319003654
$ws.Cells.Item(5, 10).Value = 1.042072747288372
$ws.Cells.Item(5, 11).Value = 1.046983939212209
$ws.Cells.Item(5, 12).Value = 1.038807330049264
$ws.Cells.Item(5, 13).Value = 1.054606727682174
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037685505543327
$ws.Cells.Item(6, 4).Value = 1.044604065931539
$ws.Cells.Item(6, 5).Value = 1.036409728095722
$ws.Cells.Item(6, 6).Value = 1.052253604378659
$ws.Cells.Item(6, 9).Value = 1.04107362252708
$ws.Cells.Item(6, 10).Value = 1.042096067691757
$ws.Cells.Item(6, 11).Value = 1.047006532597862
$ws.Cells.Item(6, 12).Value = 1.038832201393719
$ws.Cells.Item(6, 13).Value = 1.054637696630247
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037412884126627
$ws.Cells.Item(7, 4).Value = 1.044388905147593
$ws.Cells.Item(7, 5).Value = 1.036178370347967
$ws.Cells.Item(7, 6).Value = 1.051981236431108
$ws.Cells.Item(7, 9).Value = 1.041001796535876
$ws.Cells.Item(7, 10).Value = 1.041935700252936
$ws.Cells.Item(7, 11).Value = 1.046851153708946
$ws.Cells.Item(7, 12).Value = 1.038661193105571
$ws.Cells.Item(7, 13).Value = 1.054424767384248
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.036274405919202
$ws.Cells.Item(8, 4).Value = 1.043490394030603
$ws.Cells.Item(8, 5).Value = 1.035212989945838
$ws.Cells.Item(8, 6).Value = 1.050844573478876
$ws.Cells.Item(8, 9).Value = 1.040699333320647
$ws.Cells.Item(8, 10).Value = 1.041265101317679
$ws.Cells.Item(8, 11).Value = 1.046201143277331
$ws.Cells.Item(8, 12).Value = 1.037946708582463
$ws.Cells.Item(8, 13).Value = 1.053535243909376
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.034272259221552
$ws.Cells.Item(9, 4).Value = 1.04191032564214
$ws.Cells.Item(9, 5).Value = 1.03351822497897
$ws.Cells.Item(9, 6).Value = 1.048848495713122
$ws.Cells.Item(9, 9).Value = 1.040157952624649
$ws.Cells.Item(9, 10).Value = 1.040082409325875
$ws.Cells.Item(9, 11).Value = 1.04505374295752
$ws.Cells.Item(9, 12).Value = 1.036688919936835
$ws.Cells.Item(9, 13).Value = 1.0519697219371
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03294026907444
$ws.Cells.Item(10, 4).Value = 1.040859212791867
$ws.Cells.Item(10, 5).Value = 1.032392757582808
$ws.Cells.Item(10, 6).Value = 1.047522500707305
$ws.Cells.Item(10, 9).Value = 1.039791428730079
$ws.Cells.Item(10, 10).Value = 1.039293339574344
$ws.Cells.Item(10, 11).Value = 1.044287539338563
$ws.Cells.Item(10, 12).Value = 1.035851294112001
$ws.Cells.Item(10, 13).Value = 1.050927429659646
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032364178774991
$ws.Cells.Item(11, 4).Value = 1.040404627671755
$ws.Cells.Item(11, 5).Value = 1.031906476908553
$ws.Cells.Item(11, 6).Value = 1.046949472897795
$ws.Cells.Item(11, 9).Value = 1.039631400303726
$ws.Cells.Item(11, 10).Value = 1.038951534311405
$ws.Cells.Item(11, 11).Value = 1.043955478977021
$ws.Cells.Item(11, 12).Value = 1.035488822469964
$ws.Cells.Item(11, 13).Value = 1.050476453974149
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032150295317041
$ws.Cells.Item(12, 4).Value = 1.040235859381373
$ws.Cells.Item(12, 5).Value = 1.031726010701301
$ws.Cells.Item(12, 6).Value = 1.046736797438588
$ws.Cells.Item(12, 9).Value = 1.039571760875345
$ws.Cells.Item(12, 10).Value = 1.038824553849997
$ws.Cells.Item(12, 11).Value = 1.043832094769757
$ws.Cells.Item(12, 12).Value = 1.035354219590471
$ws.Cells.Item(12, 13).Value = 1.050308994683799
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032196169400141
$ws.Cells.Item(13, 4).Value = 1.040272056885355
$ws.Cells.Item(13, 5).Value = 1.031764714043583
$ws.Cells.Item(13, 6).Value = 1.04678240920795
$ws.Cells.Item(13, 9).Value = 1.039584562671626
$ws.Cells.Item(13, 10).Value = 1.038851792420162
$ws.Cells.Item(13, 11).Value = 1.043858562991032
$ws.Cells.Item(13, 12).Value = 1.035383090750169
$ws.Cells.Item(13, 13).Value = 1.050344912840161
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032346497009902
$ws.Cells.Item(14, 4).Value = 1.040390675491325
$ws.Cells.Item(14, 5).Value = 1.031891556239931
$ws.Cells.Item(14, 6).Value = 1.046931889546962
$ws.Cells.Item(14, 9).Value = 1.039626474519612
$ws.Cells.Item(14, 10).Value = 1.038941038442817
$ws.Cells.Item(14, 11).Value = 1.04394528085376
$ws.Cells.Item(14, 12).Value = 1.035477695427849
$ws.Cells.Item(14, 13).Value = 1.050462610643783
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032439132387908
$ws.Cells.Item(15, 4).Value = 1.040463771624831
$ws.Cells.Item(15, 5).Value = 1.031969729167744
$ws.Cells.Item(15, 6).Value = 1.047024012239965
$ws.Cells.Item(15, 9).Value = 1.0396522716115
$ws.Cells.Item(15, 10).Value = 1.038996023411599
$ws.Cells.Item(15, 11).Value = 1.043998705035564
$ws.Cells.Item(15, 12).Value = 1.035535989198768
$ws.Cells.Item(15, 13).Value = 1.050535135234495
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03297851653366
$ws.Cells.Item(16, 4).Value = 1.04088939391795
$ws.Cells.Item(16, 5).Value = 1.032425052785473
$ws.Cells.Item(16, 6).Value = 1.047560554779647
$ws.Cells.Item(16, 9).Value = 1.039802021517811
$ws.Cells.Item(16, 10).Value = 1.039316021366594
$ws.Cells.Item(16, 11).Value = 1.044309571104813
$ws.Cells.Item(16, 12).Value = 1.035875355030293
$ws.Cells.Item(16, 13).Value = 1.050957366802308
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03331703840869
$ws.Cells.Item(17, 4).Value = 1.041156524959712
$ws.Cells.Item(17, 5).Value = 1.0327109486945
$ws.Cells.Item(17, 6).Value = 1.047897419354643
$ws.Cells.Item(17, 9).Value = 1.039895602489991
$ws.Cells.Item(17, 10).Value = 1.039516712902177
$ws.Cells.Item(17, 11).Value = 1.044504492638904
$ws.Cells.Item(17, 12).Value = 1.036088291634855
$ws.Cells.Item(17, 13).Value = 1.051222314610032
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033514557061471
$ws.Cells.Item(18, 4).Value = 1.04131239128647
$ws.Cells.Item(18, 5).Value = 1.032877808535118
$ws.Cells.Item(18, 6).Value = 1.048094016324164
$ws.Cells.Item(18, 9).Value = 1.039950059164907
$ws.Cells.Item(18, 10).Value = 1.039633759981075
$ws.Cells.Item(18, 11).Value = 1.044618159091394
$ws.Cells.Item(18, 12).Value = 1.036212515659937
$ws.Cells.Item(18, 13).Value = 1.051376887248904
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03358191672635
$ws.Cells.Item(19, 4).Value = 1.041365546667061
$ws.Cells.Item(19, 5).Value = 1.032934720608604
$ws.Cells.Item(19, 6).Value = 1.048161069339286
$ws.Cells.Item(19, 9).Value = 1.039968605831061
$ws.Cells.Item(19, 10).Value = 1.039673667786167
$ws.Cells.Item(19, 11).Value = 1.044656911630708
$ws.Cells.Item(19, 12).Value = 1.036254876484652
$ws.Cells.Item(19, 13).Value = 1.051429598086596
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03328071152604
$ws.Cells.Item(20, 4).Value = 1.041127858800845
$ws.Cells.Item(20, 5).Value = 1.032680264248443
$ws.Cells.Item(20, 6).Value = 1.047861265628302
$ws.Cells.Item(20, 9).Value = 1.039885575322236
$ws.Cells.Item(20, 10).Value = 1.039495181925265
$ws.Cells.Item(20, 11).Value = 1.044483582288247
$ws.Cells.Item(20, 12).Value = 1.036065443297878
$ws.Cells.Item(20, 13).Value = 1.051193884786678
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032302226444815
$ws.Cells.Item(21, 4).Value = 1.040355742912228
$ws.Cells.Item(21, 5).Value = 1.031854199941622
$ws.Cells.Item(21, 6).Value = 1.046887866534261
$ws.Cells.Item(21, 9).Value = 1.039614137977768
$ws.Cells.Item(21, 10).Value = 1.038914758222056
$ws.Cells.Item(21, 11).Value = 1.043919745759685
$ws.Cells.Item(21, 12).Value = 1.035449835725399
$ws.Cells.Item(21, 13).Value = 1.050427950094546
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031687605642651
$ws.Cells.Item(22, 4).Value = 1.039870774832715
$ws.Cells.Item(22, 5).Value = 1.031335747586074
$ws.Cells.Item(22, 6).Value = 1.046276851292288
$ws.Cells.Item(22, 9).Value = 1.039442330820628
$ws.Cells.Item(22, 10).Value = 1.038549714808012
$ws.Cells.Item(22, 11).Value = 1.043564995788027
$ws.Cells.Item(22, 12).Value = 1.035062983265976
$ws.Cells.Item(22, 13).Value = 1.049946684864848
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032013371201166
$ws.Cells.Item(23, 4).Value = 1.040127818407934
$ws.Cells.Item(23, 5).Value = 1.031610500542396
$ws.Cells.Item(23, 6).Value = 1.046600666701554
$ws.Cells.Item(23, 9).Value = 1.039533517208406
$ws.Cells.Item(23, 10).Value = 1.03874324106725
$ws.Cells.Item(23, 11).Value = 1.043753078134826
$ws.Cells.Item(23, 12).Value = 1.03526804126252
$ws.Cells.Item(23, 13).Value = 1.050201782938829
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033297125881278
$ws.Cells.Item(24, 4).Value = 1.041140811639603
$ws.Cells.Item(24, 5).Value = 1.032694128915125
$ws.Cells.Item(24, 6).Value = 1.047877601602702
$ws.Cells.Item(24, 9).Value = 1.039890106561836
$ws.Cells.Item(24, 10).Value = 1.039504910886932
$ws.Cells.Item(24, 11).Value = 1.044493030862231
$ws.Cells.Item(24, 12).Value = 1.036075767410948
$ws.Cells.Item(24, 13).Value = 1.051206730898116
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03478937950228
$ws.Cells.Item(25, 4).Value = 1.042318419033517
$ws.Cells.Item(25, 5).Value = 1.033955598429503
$ws.Cells.Item(25, 6).Value = 1.049363704881868
$ws.Cells.Item(25, 9).Value = 1.040298903491952
$ws.Cells.Item(25, 10).Value = 1.040388276015139
$ws.Cells.Item(25, 11).Value = 1.045350602710318
$ws.Cells.Item(25, 12).Value = 1.037013935447001
$ws.Cells.Item(25, 13).Value = 1.05237420865781
